# "Loan RBI, Variable Instalments"
#
# The "Repayment schedule" sheet gets a new (blank) column inserted before
# the existing "Late" column, pushing Late / Outstanding(heading) / Outstanding
# one column to the right (N->O, O->P, P->Q) and widening the used range
# from A1:P14 to A1:Q14. The new column N inherits the column width of the
# column immediately to its left (M).
#
# The workbook is also left with the "Repayment schedule" tab active/selected
# (instead of "Input"), with cell R5 selected on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N - shifts old N/O/P (Late / heading / Outstanding)
# one column to the right, and extends the sheet dimension to A1:Q14.
$ws.Columns("N").Insert()

# Excel's column insert carries the left-neighbour's width onto the new column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active/selected sheet (was "Input"), with R5 selected.
$ws.Activate() | Out-Null
$ws.Range("R5").Select() | Out-Null
